$d = $word.ActiveDocument

# --- Paragraph 1: "Package Installation guid" (title) ---
# Merge the two existing runs ("Package Installation " + "guid") into a
# single run that keeps the first run's rsid, drop the spell-check
# proofErr markers around "guid", and append a new run containing "e"
# (completing the word "guide") with matching bold formatting.
$p1 = $d.Paragraphs(1).Range
$p1xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00245F4C" w:rsidRPr="00BF7838" w:rsidRDefault="00515695" w:rsidP="00515695"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr><w:r w:rsidRPr="00BF7838"><w:rPr><w:b/></w:rPr><w:t>Package Installation guid</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>e</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p1.InsertXML($p1xml)

# --- Paragraph 2: "1. Select the "Manage Package" in Uipath Studio" ---
# Merge the three existing runs into a single run (drop the proofErr
# markers that wrapped "Uipath").
$p2 = $d.Paragraphs(2).Range
$p2xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00515695" w:rsidRPr="00515695" w:rsidRDefault="00515695" w:rsidP="00515695"><w:r><w:t>1. Select the &quot;Manage Package&quot; in Uipath Studio</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p2.InsertXML($p2xml)
